# Adds rows 60-62 to the chat registry sheet (restores entries removed by an obsolete test-script cleanup).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$a60 = @'
PARTE POLICIAL N.º 093-2025
UNIDAD: Policía Nacional del Ecuador – Unidad de Investigación de Transparencia Institucional (UITI)
FECHA: 29 de junio de 2025
HORA DEL REPORTE: 17h50
LUGAR DE LOS HECHOS: GAD Parroquial de San Ramón, Provincia de Imbabura
TIPO DE INCIDENTE: Investigación preliminar por presunto peculado
DENUNCIANTE: Ciudadano anónimo (vía buzón virtual de la Contraloría General del Estado)
FUNCIONARIO ACTUANTE: Sargento Primero Paulina Enríquez C.
HECHOS:
Mediante oficio N.º CGE-OF-0032-2025 remitido por la Contraloría General del Estado, se dispone iniciar investigación preliminar por presunto peculado relacionado con la administración de fondos públicos destinados al proyecto “Mejoramiento del sistema de riego El Cebadal”, ejecutado por el GAD Parroquial de San Ramón durante el segundo semestre de 2024.
La denuncia anónima sostiene que el monto asignado, por un total de **$186,000 dólares**, habría sido desembolsado sin que las obras estén concluidas ni verificados los avances técnicos. Se señala al **presidente del GAD parroquial, señor Galo Yánez Salinas**, como la autoridad responsable de suscribir las actas de recepción de obra, aunque no se cuenta aún con pruebas concluyentes sobre su participación directa en un acto doloso.
Durante el operativo inicial, se recabó documentación financiera, contratos, y copias de actas entregadas por la entidad. Asimismo, se entrevistó a tres funcionarios técnicos que manifestaron que hubo retrasos en la ejecución, pero que no podían confirmar un mal uso de los recursos.
La Fiscalía Provincial de Imbabura fue notificada y se encuentra valorando el inicio formal de una instrucción fiscal. El ciudadano investigado fue notificado de sus derechos conforme lo establece el artículo 77 de la Constitución y el COIP, y ha manifestado estar dispuesto a colaborar con la investigación.
**OBSERVACIONES ADICIONALES:** Se aclara que hasta el momento **no existe orden de detención ni formulación de cargos**, y que el presente parte se emite únicamente con fines informativos y de coordinación interinstitucional.
FIRMA DEL FUNCIONARIO:  Sargento Primero Paulina Enríquez C.  Policía Nacional del Ecuador
'@
$ws.Cells.Item(60, 1).Value = $a60

$b60 = @'
📅 **FECHA Y HORA:**
01 de July de 2025, 15:56 horas
⚖️ **RAZÓN DE LA SENTENCIA:**
De acuerdo con el Código Orgánico Integral Penal (COIP) ecuatoriano, en el caso presentado se ha llevado a cabo una investigación preliminar por presunto peculado relacionado con la administración de fondos públicos destinados al proyecto "Mejoramiento del sistema de riego El Cebadal" ejecutado por el GAD Parroquial de San Ramón. La denuncia anónima indica que el monto asignado de $186,000 dólares habría sido desembolsado sin que las obras estén concluidas ni verificados los avances técnicos, señalando al presidente del GAD parroquial, señor Galo Yánez Salinas, como la autoridad responsable. A pesar de no contar con pruebas concluyentes sobre su participación directa en un acto doloso, se ha iniciado una investigación por parte de la Fiscalía Provincial de Imbabura.
🏛️ **VEREDICTO:**
En virtud de lo establecido en el COIP, y considerando que hasta el momento no existe orden de detención ni formulación de cargos contra el ciudadano investigado, se declara al señor Galo Yánez Salinas como **INOCENTE** en este proceso.
🏢 **LUGAR DE RECLUSIÓN:**
Libertad inmediata del procesado, señor Galo Yánez Salinas.
📋 **CONCLUSIÓN:**
En atención a lo dispuesto por el Código Orgánico Integral Penal ecuatoriano, se determina que no existen elementos suficientes para considerar al señor Galo Yánez Salinas culpable de los cargos de peculado en el presente caso. Se ordena su libertad inmediata y se le recuerdan sus derechos conforme al artículo 77 de la Constitución y el COIP. Notifíquese y cúmplase.
NOTIFÍQUESE Y CÚMPLASE.
📚 **Información de consulta:**
- Documento: Código Orgánico Integral Penal
- Especialidad: Derecho Penal
- Fragmentos consultados: 8
- Relevancia máxima: 0.751
- Relevancia mínima: 0.736
'@
$ws.Cells.Item(60, 2).Value = $b60

$ws.Cells.Item(60, 3).Value = "'46"

$a61 = @'
[SIN CONTEXTO DOC] Hola
'@
$ws.Cells.Item(61, 1).Value = $a61

$b61 = @'
📅 **FECHA Y HORA:**
01 de July de 2025, 16:07 horas
⚖️ **RAZÓN DE LA SENTENCIA:**
No se cuenta con el marco legal específico en el documento proporcionado para este caso. Sin embargo, basándome en principios generales del derecho ecuatoriano:
🏛️ **VEREDICTO:**
SUSPENDIDO - Falta de marco legal específico en el documento proporcionado
🏢 **LUGAR DE RECLUSIÓN:**
No aplicable hasta contar con legislación específica
📋 **CONCLUSIÓN:**
Esta sentencia provisional se basa en conocimiento general jurídico y NO en el documento específico proporcionado.
Se requiere consultar la legislación correspondiente y asesoramiento legal especializado.
CASO PRESENTADO:
Hola
⚠️ ADVERTENCIA JUDICIAL: Sentencia no definitiva por falta de marco legal específico.
'@
$ws.Cells.Item(61, 2).Value = $b61

$ws.Cells.Item(61, 3).Value = "'41"

$a62 = @'
Tengo 12 años y me obligan a trabajar
'@
$ws.Cells.Item(62, 1).Value = $a62

$b62 = @'
📅 **FECHA Y HORA:**
14 de noviembre de 2021, 15:30.
⚖️ **RAZÓN DE LA SENTENCIA:**
El caso presentado versa sobre la explotación laboral de una persona menor de quince años, lo cual constituye una grave violación a los derechos humanos y a la normativa ecuatoriana. Según el Art. 105 del None, se establece que "La persona que someta a otra a trabajos forzados u otras formas de explotación o servicios laborales, dentro o fuera del país, será sancionada con pena privativa de libertad de diez a trece años". Asimismo, el Art. 69 del Código de la Niñez y Adolescencia dispone que se debe proteger a los menores de cualquier forma de explotación laboral.
🏛️ **VEREDICTO:**
CULPABLE. En virtud de lo establecido en el Art. 105 del None y en concordancia con el Art. 69 del Código de la Niñez y Adolescencia, se declara al acusado como culpable de someter a una persona menor de quince años a trabajos forzados.
🏢 **LUGAR DE RECLUSIÓN:**
El acusado será recluido en el centro penitenciario correspondiente, de acuerdo a la gravedad del delito cometido, por un periodo de diez a trece años, conforme a lo establecido en el Art. 105 del None.
📋 **CONCLUSIÓN:**
En atención a lo expuesto y a las disposiciones legales pertinentes, se determina que la explotación laboral de menores de edad es un delito grave que atenta contra los derechos fundamentales de los niños y adolescentes. Por tanto, se condena al acusado a una pena privativa de libertad de diez a trece años, de conformidad con lo establecido en el Art. 105 del None. Notifíquese y cúmplase.
NOTIFÍQUESE Y CÚMPLASE.
📚 **Información de consulta:**
- Documento: None
- Especialidad: None
- Fragmentos consultados: 8
- Relevancia máxima: 0.587
- Relevancia mínima: 0.477
'@
$ws.Cells.Item(62, 2).Value = $b62

$ws.Cells.Item(62, 3).Value = "'45"

